$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 2064.03
$ws.Range("H3").Value = 0.0004844890820385362
$ws.Range("I4").Value = 2064.03
$ws.Range("H5").Value = 0.0004842700929734548
$ws.Range("I5").Value = 0.999548
